# "update the cycle model" - rename the crop_type column/value and fix the
# dli_u cell so it holds a real number instead of the placeholder text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crop")

# C2 must be written before C1 so the shared-string table ends up with
# AmaranthRed appended ahead of type_crop (matches the saved workbook order).
$ws.Range("C2").Value = "AmaranthRed"      # crop name: BayamRed -> AmaranthRed
$ws.Range("C1").Value = "type_crop"        # header: crop_type -> type_crop

$ws.Range("O2").Value = 12                 # dli_u: "unknown" -> 12

# Move the selection to O3, matching the saved view state.
$ws.Range("O3").Select()
